$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17 is the "Case File" rule. Flip "Start a Workflow Process?" (column D)
# from true to false, and clear out the now-irrelevant process name (E),
# priority (F) and due-date (G) columns that only apply when a workflow starts.
$ws.Cells.Item(17, 4).Value = "'false"
$ws.Cells.Item(17, 4).NumberFormat = '"TRUE";"TRUE";"FALSE"'

$ws.Cells.Item(17, 5).ClearContents()
$ws.Cells.Item(17, 6).ClearContents()
$ws.Cells.Item(17, 7).ClearContents()

# Reflect the author's final cursor position in the sheet view.
$ws.Range("G18").Select()
